# Microns per pixel is different for the two microscopes, so this needs to be
# set explicitly. Add a new settings line (with explanatory comments) just
# below the micromanager flag on the "pos5_B" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "pos5_B" - the active/selected sheet

# Insert a new blank row above row 36 (i.e. right after the existing
# "%p1.micromanager=1;" line in row 35). Everything from old row 36 onward
# shifts down by one; the new row inherits row 35's formatting.
$ws.Rows.Item(36).Insert()

# New row 36, column A: the actual setting.
$ws.Range("A36").Value = "%p.micronsPerPixel = 0.0431;"

# Row 35, column C: a short comment next to the pre-existing micromanager line.
$ws.Range("C35").Value = "% Only for Tans2 microscope"
$ws.Range("C35").Font.Name = "Arial"
$ws.Range("C35").Font.Size = 8
$ws.Range("C35").Font.Bold = $true
$ws.Range("C35").Font.Color = 16711680

# New row 36, column C: explanatory comment next to the new setting.
$ws.Range("C36").Value = "% Only for Tans2 microscope - Note that MetaData dir needs to be in root pos dir."
$ws.Range("C36").Font.Name = "Arial"
$ws.Range("C36").Font.Size = 8
$ws.Range("C36").Font.Bold = $true
$ws.Range("C36").Font.Color = 16711680

# Leave the cursor where the edit happened.
$ws.Range("B40").Select()
